$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "lines" table (line1..line6) is being extended with two more line
# entries (line7, line8). Insert two new rows right above the first
# "extr" row (row 8) so the extr1..extr8 block shifts down intact.
$ws.Rows("8:9").Insert()

# Match the formatting of the surrounding data rows (Insert() otherwise
# leaves the new rows with a slightly different inherited style).
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E9").PasteSpecial(-4122)

# Fill in the two new "line" rows.
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Renumber the "name" index column for the shifted extr1..extr8 rows
# (now rows 10..17) so it keeps counting up without a gap.
for ($r = 10; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Flip in_service for a few of the contingencies.
$ws.Cells.Item(10, 5).Value = $true   # extr1
$ws.Cells.Item(12, 5).Value = $true   # extr3
$ws.Cells.Item(13, 5).Value = $true   # extr4
$ws.Cells.Item(14, 5).Value = $false  # extr5
